$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "243.40"
    "D3"  = "23.11"
    "D4"  = "5.398"
    "D6"  = "3.394"
    "D8"  = "0.9110"
    "D9"  = "0.1415"
    "D10" = "0.07437"
    "D11" = "0.03317"
    "D12" = "0.03066"
    "D13" = "0.09327"
    "D14" = "3.960"
    "D15" = "0.001578"
    "D16" = "0.04810"
    "D17" = "0.0005942"
    "D18" = "0.006142"
    "D20" = "0.004428"
    "D21" = "0.0009839"
    "D22" = "0.00007802"
    "D23" = "3.626"
    "D25" = "2.149"
    "D40" = "0.03874"
    "D41" = "0.006198"
    "D42" = "0.1067"
    "D43" = "0.002801"
    "D44" = "0.006614"
    "D45" = "0.00005186"
    "D47" = "0.0005802"
    "D48" = "0.8249"
    "D49" = "0.002261"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
